$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing IPC values (column B) for rows 2-21
$ws.Range("B2").Value = 2.609749512915265
$ws.Range("B3").Value = 2.599760133123742
$ws.Range("B4").Value = 2.381642940183032
$ws.Range("B5").Value = 2.376115025220395
$ws.Range("B6").Value = 2.465060601522592
$ws.Range("B7").Value = 2.466574831877659
$ws.Range("B8").Value = 2.086868822678343
$ws.Range("B9").Value = 2.083833551784116
$ws.Range("B10").Value = 1.904062509730511
$ws.Range("B11").Value = 1.900953431346346
$ws.Range("B12").Value = 1.828997671118578
$ws.Range("B13").Value = 1.83194526512898
$ws.Range("B14").Value = 2.331243515450478
$ws.Range("B15").Value = 2.332827623620034
$ws.Range("B16").Value = 2.079387492158941
$ws.Range("B17").Value = 1.878180692744194
$ws.Range("B18").Value = 1.783929705295721
$ws.Range("B19").Value = 1.976959778621488
$ws.Range("B20").Value = 1.856864716022876
$ws.Range("B21").Value = 1.795273704746372

# Add new row 22 (A22 = 24, B22 = new IPC value) - copy A21 formatting to A22
$ws.Range("A21").Copy($ws.Range("A22"))
$ws.Range("A22").Value = 24
$ws.Range("B22").Value = 2.908958354580876
